# Section 1 edit:
#   1. Remove page 2  -> slide 3's "V2"/"V2" placeholders become "V1.a" (title)
#      and an emptied content placeholder (the old "page 2" content is gone).
#   2. Add page 1.a   -> slide 1's subtitle gets a new paragraph "v1.a"
#      underneath the existing "Steven" line.

$p = $ppt.ActivePresentation

# --- Slide 1: title slide ("Steven") -> add a new "v1.a" paragraph ---
$s1 = $p.Slides.Item(1)
$subtitle = $s1.Shapes.Item("Subtitle 2")
$subTr = $subtitle.TextFrame.TextRange

# Start a new paragraph after "Steven", then type "v" and "1.a" as two
# separate runs (mirrors the authored edit, which shows the text typed
# in two pieces).
[void]$subTr.InsertAfter([char]13 + "v")
[void]$subtitle.TextFrame.TextRange.InsertAfter("1.a")

# --- Slide 3: "V2" page -> retitled "V1.a", old body content removed ---
$s3 = $p.Slides.Item(3)

$title3 = $s3.Shapes.Item("Title 1")
$title3.TextFrame.TextRange.Text = "V1.a"

$content3 = $s3.Shapes.Item("Content Placeholder 2")
$contentTr = $content3.TextFrame.TextRange
$contentLen = $contentTr.Length
if ($contentLen -gt 0) {
    $contentTr.Start = 1
    [void]$contentTr.Characters(1, $contentLen).Delete()
}
